$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '61.489.73'
$ws.Range("D2").Style = "Normal"

$ws.Range("E2").Value = '  +1.07%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.442.12'
$ws.Range("D3").Style = "Normal"

$ws.Range("E3").Value = '  +1.95%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("D4").Style = "Normal"

$ws.Range("E4").Value = '  +0.00%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '576.83'
$ws.Range("D5").Style = "Normal"

$ws.Range("E5").Value = '  +1.43%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '144.63'
$ws.Range("D6").Style = "Normal"

$ws.Range("E6").Value = '  +6.33%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '3.444.13'
$ws.Range("D7").Style = "Normal"

$ws.Range("E7").Value = '  +2.05%  '

$ws.Range("E8").Value = '  +0.04%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.476'
$ws.Range("D9").Style = "Normal"

$ws.Range("E9").Value = '  +1.90%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '7.63'
$ws.Range("D10").Style = "Normal"

$ws.Range("E10").Value = '  +0.33%  '

$ws.Range("E11").Value = '  +3.41%  '

$ws.Range("E12").Value = '  +1.91%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.030.01'
$ws.Range("D13").Style = "Normal"

$ws.Range("E13").Value = '  +2.00%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '28.02'
$ws.Range("D14").Style = "Normal"

$ws.Range("E14").Value = '  +8.03%  '

$ws.Range("E15").Value = '  -1.00%  '

$ws.Range("B16").Value = 'ShibaInu'

$ws.Range("C16").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.0000172'
$ws.Range("D16").Style = "Normal"

$ws.Range("E16").Value = '  +2.03%  '

$ws.Range("B17").Value = 'WrappedEther'

$ws.Range("C17").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.437.24'
$ws.Range("D17").Style = "Normal"

$ws.Range("E17").Value = '  +1.96%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '61.604.15'
$ws.Range("D18").Style = "Normal"

$ws.Range("E18").Value = '  +1.09%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.26'
$ws.Range("D19").Style = "Normal"

$ws.Range("E19").Value = '  +8.06%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '14.18'
$ws.Range("D20").Style = "Normal"

$ws.Range("E20").Value = '  +3.85%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '9.49'
$ws.Range("D21").Style = "Normal"

$ws.Range("E21").Value = '  +2.96%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '395.45'
$ws.Range("D22").Style = "Normal"

$ws.Range("E22").Value = '  +6.44%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.563'
$ws.Range("D23").Style = "Normal"

$ws.Range("E23").Value = '  +3.18%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '73.15'
$ws.Range("D24").Style = "Normal"

$ws.Range("E24").Value = '  +3.39%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.999'
$ws.Range("D25").Style = "Normal"

$ws.Range("E25").Value = '  -0.25%  '

$ws.Range("E26").Value = '  -0.25%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.0000122'
$ws.Range("D27").Style = "Normal"

$ws.Range("E27").Value = '  -0.01%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '3.588.74'
$ws.Range("D28").Style = "Normal"

$ws.Range("E28").Value = '  +2.34%  '

$ws.Range("E29").Value = '  +0.63%  '

$ws.Range("E30").Value = '  +3.85%  '

$ws.Range("E31").Value = '  +0.08%  '

$ws.Range("B32").Value = 'Fetch.AI'

$ws.Range("C32").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.48'
$ws.Range("D32").Style = "Normal"

$ws.Range("E32").Value = '  -8.14%  '

$ws.Range("B33").Value = 'InternetComputer(DFINITY)'

$ws.Range("C33").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '8.15'
$ws.Range("D33").Style = "Normal"

$ws.Range("E33").Value = '  +1.97%  '

$ws.Range("E34").Value = '  +2.11%  '

$ws.Range("E35").Value = '  -0.07%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '23.97'
$ws.Range("D36").Style = "Normal"

$ws.Range("E36").Value = '  +3.09%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.472.97'
$ws.Range("D37").Style = "Normal"

$ws.Range("E37").Value = '  +2.27%  '

$ws.Range("E38").Value = '  +3.66%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.12'
$ws.Range("D39").Style = "Normal"

$ws.Range("E39").Value = '  +0.56%  '

$ws.Range("E40").Value = '  +1.63%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '167.53'
$ws.Range("D41").Style = "Normal"

$ws.Range("E41").Value = '  +1.72%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.0780'
$ws.Range("D42").Style = "Normal"

$ws.Range("E42").Value = '  +3.20%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '27.53'
$ws.Range("D43").Style = "Normal"

$ws.Range("E43").Value = '  +10.05%  '

$ws.Range("E44").Value = '  +4.01%  '

$ws.Range("E45").Value = '  +0.05%  '

$ws.Range("E46").Value = '  +2.10%  '

$ws.Range("E47").Value = '  +1.33%  '

$ws.Range("E48").Value = '  +4.03%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.600.44'
$ws.Range("D49").Style = "Normal"

$ws.Range("E49").Value = '  +2.70%  '

$ws.Range("E50").Value = '  -0.31%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '6.91'
$ws.Range("D51").Style = "Normal"

$ws.Range("E51").Value = '  +2.46%  '
